$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column B. Target stored width (OOXML <col width>) is 20.375 chars,
# matching column D. The COM ColumnWidth setter here rounds to whole pixels
# (stored = round(ColumnWidth*7 + 5) / 7), so 20.375 itself is not a
# reachable value; 19.714286 is the input that lands on the closest
# reachable stored width (20.428571...).
$ws.Columns.Item(2).ColumnWidth = 19.714286

# Apply centered style (alignment) to the numeric/label cells in both
# In (C14:D18) and Out (B24:C27) tables so they match the rest of the table.
$centerRanges = @("C14","D14","C15","D15","C16","D16","C17","D17","C18","D18", `
                   "B24","C24","B25","C25","B26","C26","B27","C27")

foreach ($addr in $centerRanges) {
    $rng = $ws.Range($addr)
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4108
}

# Move the selection to D10
$ws.Range("D10").Select()
